$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(334, 6).Value = 195214
$ws.Cells.Item(334, 7).Value = 3437
$ws.Cells.Item(335, 6).Value = 130729
$ws.Cells.Item(335, 7).Value = 2989
$ws.Cells.Item(336, 6).Value = 101314
$ws.Cells.Item(336, 7).Value = 3347
$ws.Cells.Item(337, 6).Value = 103994
$ws.Cells.Item(337, 7).Value = 2962
$ws.Cells.Item(338, 6).Value = 226341
$ws.Cells.Item(338, 7).Value = 3169
$ws.Cells.Item(339, 6).Value = 656936
$ws.Cells.Item(339, 7).Value = 5475
$ws.Cells.Item(340, 6).Value = 383003
$ws.Cells.Item(340, 7).Value = 3299
$ws.Cells.Item(341, 6).Value = 291643
$ws.Cells.Item(341, 7).Value = 3653
$ws.Cells.Item(342, 6).Value = 179374
$ws.Cells.Item(342, 7).Value = 3060
$ws.Cells.Item(343, 6).Value = 132294
$ws.Cells.Item(343, 7).Value = 2934
$ws.Cells.Item(344, 6).Value = 135213
$ws.Cells.Item(344, 7).Value = 2471
$ws.Cells.Item(345, 6).Value = 291437
$ws.Cells.Item(345, 7).Value = 3300
$ws.Cells.Item(346, 6).Value = 669364
$ws.Cells.Item(346, 7).Value = 4783
$ws.Cells.Item(347, 6).Value = 341714
$ws.Cells.Item(347, 7).Value = 2891
$ws.Cells.Item(348, 6).Value = 232063
$ws.Cells.Item(348, 7).Value = 3243
$ws.Cells.Item(349, 6).Value = 159815
$ws.Cells.Item(349, 7).Value = 2744
$ws.Cells.Item(350, 6).Value = 126998
$ws.Cells.Item(350, 7).Value = 2777
$ws.Cells.Item(351, 6).Value = 150233
$ws.Cells.Item(351, 7).Value = 2820
$ws.Cells.Item(352, 6).Value = 306400
$ws.Cells.Item(352, 7).Value = 3534
$ws.Cells.Item(353, 6).Value = 720004
$ws.Cells.Item(353, 7).Value = 5254
$ws.Cells.Item(354, 6).Value = 308829
$ws.Cells.Item(354, 7).Value = 2826
$ws.Cells.Item(355, 6).Value = 222632
$ws.Cells.Item(356, 6).Value = 159991
$ws.Cells.Item(356, 7).Value = 2888
$ws.Cells.Item(357, 6).Value = 138321
$ws.Cells.Item(357, 7).Value = 3027
$ws.Cells.Item(358, 6).Value = 157206
$ws.Cells.Item(358, 7).Value = 2599
$ws.Cells.Item(359, 6).Value = 320737
$ws.Cells.Item(359, 7).Value = 3338
$ws.Cells.Item(360, 6).Value = 743526
$ws.Cells.Item(361, 6).Value = 331526
$ws.Cells.Item(362, 6).Value = 227056
$ws.Cells.Item(362, 7).Value = 3139
$ws.Cells.Item(363, 6).Value = 187459
$ws.Cells.Item(363, 7).Value = 2748
$ws.Cells.Item(364, 6).Value = 166863
$ws.Cells.Item(364, 7).Value = 2449
$ws.Cells.Item(365, 6).Value = 179525
$ws.Cells.Item(365, 7).Value = 2344
$ws.Cells.Item(366, 6).Value = 335824
$ws.Cells.Item(366, 7).Value = 2814
$ws.Cells.Item(367, 6).Value = 753389
$ws.Cells.Item(367, 7).Value = 3852
$ws.Cells.Item(368, 6).Value = 342891
$ws.Cells.Item(369, 6).Value = 231849
$ws.Cells.Item(369, 7).Value = 2559
$ws.Cells.Item(370, 6).Value = 180627
$ws.Cells.Item(370, 7).Value = 2014
$ws.Cells.Item(371, 6).Value = 157189
$ws.Cells.Item(371, 7).Value = 1923
$ws.Cells.Item(372, 6).Value = 174422
$ws.Cells.Item(372, 7).Value = 1807
$ws.Cells.Item(373, 6).Value = 339043
$ws.Cells.Item(373, 7).Value = 2299
$ws.Cells.Item(374, 6).Value = 732321
$ws.Cells.Item(374, 7).Value = 3286
$ws.Cells.Item(375, 6).Value = 334316
$ws.Cells.Item(375, 7).Value = 1791
$ws.Cells.Item(376, 6).Value = 214856
$ws.Cells.Item(376, 7).Value = 2144
